$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-06-21 Saturday" "2025-06-22 Sunday"
Replace-Text "337×4=1348" "818×5=4090"
Replace-Text "991×9=8919" "773×5=3865"
Replace-Text "430×9=3870" "354×6=2124"
Replace-Text "378×5=1890" "304×9=2736"
Replace-Text "946×8=7568" "531×3=1593"
Replace-Text "893×8=7144" "855×3=2565"
Replace-Text "691×8=5528" "577×8=4616"
Replace-Text "409×7=2863" "585×2=1170"
Replace-Text "726×4=2904" "976×4=3904"
Replace-Text "825×4=3300" "313×8=2504"
Replace-Text "301×6=1806" "908×2=1816"
Replace-Text "170×8=1360" "995×8=7960"
Replace-Text "894×9=8046" "594×5=2970"
Replace-Text "980×4=3920" "975×2=1950"
Replace-Text "525×9=4725" "796×3=2388"
Replace-Text "435×8=3480" "280×7=1960"
Replace-Text "284×4=1136" "134×3=402"
Replace-Text "173×3=519" "794×3=2382"
Replace-Text "147×4=588" "538×7=3766"
Replace-Text "449×5=2245" "355×9=3195"
Replace-Text "387×7=2709" "728×6=4368"
Replace-Text "923×5=4615" "980×8=7840"
Replace-Text "274×6=1644" "588×9=5292"
Replace-Text "516×6=3096" "964×7=6748"
Replace-Text "244×6=1464" "219×3=657"
